$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with freshly scraped
# coinranking.com values. Values are stored as literal text (matching the
# workbook's existing inline-string cells), so each is written with a
# leading apostrophe to stop Excel from reinterpreting the numeric-looking
# / percent-looking text as a Number value.

$ws.Range("D2").Value = "'308.95"
$ws.Range("E2").Value = "'0.33%"
$ws.Range("D3").Value = "'41.12"
$ws.Range("E3").Value = "'0.31%"
$ws.Range("D4").Value = "'5.183"
$ws.Range("E4").Value = "'1.47%"
$ws.Range("D5").Value = "'0.07686"
$ws.Range("E5").Value = "'0.76%"
$ws.Range("D6").Value = "'1.651"
$ws.Range("E6").Value = "'2.92%"
$ws.Range("D7").Value = "'0.9144"
$ws.Range("E7").Value = "'1.13%"
$ws.Range("D9").Value = "'0.1242"
$ws.Range("E9").Value = "'11.27%"
$ws.Range("D10").Value = "'0.1822"
$ws.Range("E10").Value = "'2.47%"
$ws.Range("D11").Value = "'0.09235"
$ws.Range("E11").Value = "'1.10%"
$ws.Range("D12").Value = "'0.04218"
$ws.Range("E12").Value = "'0.07%"
$ws.Range("E13").Value = "'-0.02%"
$ws.Range("D14").Value = "'0.001256"
$ws.Range("E14").Value = "'-0.17%"
$ws.Range("D15").Value = "'0.005754"
$ws.Range("E15").Value = "'0.85%"
$ws.Range("E16").Value = "'1,903.15%"
$ws.Range("E17").Value = "'-0.07%"
$ws.Range("E18").Value = "'1.64%"
$ws.Range("D20").Value = "'7.399"
$ws.Range("E20").Value = "'11.33%"
$ws.Range("D21").Value = "'0.1382"
$ws.Range("E21").Value = "'1.23%"
$ws.Range("D22").Value = "'0.2819"
$ws.Range("E22").Value = "'0.61%"
$ws.Range("D23").Value = "'0.04021"
$ws.Range("E23").Value = "'-1.24%"
$ws.Range("E24").Value = "'1.57%"
$ws.Range("D25").Value = "'0.004090"
$ws.Range("E25").Value = "'-0.45%"
$ws.Range("D26").Value = "'0.0001300"
$ws.Range("E26").Value = "'-0.09%"
$ws.Range("D38").Value = "'0.02559"
$ws.Range("E38").Value = "'7.05%"
$ws.Range("D39").Value = "'0.05321"
$ws.Range("E39").Value = "'2.42%"
$ws.Range("D40").Value = "'0.007838"
$ws.Range("E40").Value = "'0.77%"
$ws.Range("D41").Value = "'0.1315"
$ws.Range("E41").Value = "'0.92%"
$ws.Range("D42").Value = "'0.006663"
$ws.Range("E42").Value = "'-5.56%"
$ws.Range("E43").Value = "'-4.70%"
$ws.Range("D44").Value = "'0.007997"
$ws.Range("E44").Value = "'0.70%"
$ws.Range("D45").Value = "'0.3071"
$ws.Range("E45").Value = "'-0.45%"
$ws.Range("D46").Value = "'0.00006726"
$ws.Range("E46").Value = "'-4.09%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.08%"
$ws.Range("D48").Value = "'0.2629"
$ws.Range("E48").Value = "'727.46%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.08%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.08%"
